$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1424.2858
$ws.Range("I19").Value = 992.5
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 992.5
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -817.5
$ws.Range("N19").Value = -2350
$ws.Range("H33").Value = 45635.684
$ws.Range("I33").Value = 52831.684
$ws.Range("K33").Value = 52831.684
$ws.Range("M33").Value = -52602.684
$ws.Range("H70").Value = 1471.7142
$ws.Range("I70").Value = 1434
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 4302
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -4032
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1471.7142
$ws.Range("I73").Value = 1434
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 4302
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -3366
$ws.Range("N73").Value = -6372
$ws.Range("H138").Value = 4025.23
$ws.Range("I138").Value = 897.1429
$ws.Range("J138").Value = 4625.137
$ws.Range("K138").Value = 2691.4287
$ws.Range("L138").Value = 13875.411
$ws.Range("M138").Value = 2448.5713
$ws.Range("N138").Value = -24155.411
$ws.Range("H141").Value = 2798.3333
$ws.Range("I141").Value = 2082.6924
$ws.Range("K141").Value = 6248.0772
$ws.Range("M141").Value = -1068.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2489.1428
$ws.Range("I61").Value = 2570.8333
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 2570.8333
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -2358.8333
$ws.Range("N61").Value = -2423
$ws.Range("H128").Value = 49490
$ws.Range("J128").Value = 49490
$ws.Range("L128").Value = 49490
$ws.Range("N128").Value = -59450
$ws.Range("H136").Value = 2489.1428
$ws.Range("I136").Value = 2570.8333
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 7712.499899999999
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -5162.499899999999
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2185.7334
$ws.Range("I86").Value = 1732.1666
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 1732.1666
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -609.1666
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 2185.7334
$ws.Range("I89").Value = 1732.1666
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 8660.833
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -3044.833000000001
$ws.Range("N89").Value = -31232
$ws.Range("H132").Value = 109650
$ws.Range("J132").Value = 109650
$ws.Range("L132").Value = 109650
$ws.Range("N132").Value = -119770
$ws.Range("H137").Value = 55338.46
$ws.Range("J137").Value = 55338.46
$ws.Range("L137").Value = 55338.46
$ws.Range("N137").Value = -65538.45999999999
$ws.Range("H138").Value = 40130
$ws.Range("J138").Value = 40130
$ws.Range("L138").Value = 40130
$ws.Range("N138").Value = -50410
$ws.Range("H140").Value = 73972.375
$ws.Range("J140").Value = 73972.375
$ws.Range("L140").Value = 73972.375
$ws.Range("N140").Value = -84332.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H69").Value = 1360
$ws.Range("J69").Value = 1600
$ws.Range("L69").Value = 4800
$ws.Range("N69").Value = -6422
$ws.Range("H72").Value = 1360
$ws.Range("J72").Value = 1600
$ws.Range("L72").Value = 14400
$ws.Range("N72").Value = -22512
$ws.Range("H82").Value = 11900
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 13550
$ws.Range("K82").Value = 6000
$ws.Range("L82").Value = 40650
$ws.Range("M82").Value = -5594
$ws.Range("N82").Value = -41462
$ws.Range("H85").Value = 11900
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 13550
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 40650
$ws.Range("M85").Value = -4596
$ws.Range("N85").Value = -43458
$ws.Range("H108").Value = 1470.4
$ws.Range("I108").Value = 1470.4
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 4411.200000000001
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -1531.200000000001
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 546.26
$ws.Range("I113").Value = 504.5625
$ws.Range("J113").Value = 565.8823
$ws.Range("K113").Value = 1513.6875
$ws.Range("L113").Value = 1697.6469
$ws.Range("M113").Value = 656.3125
$ws.Range("N113").Value = -6037.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1628.5714
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -6140
$ws.Range("H134").Value = 18790.4
$ws.Range("J134").Value = 18790.4
$ws.Range("L134").Value = 56371.2
$ws.Range("N134").Value = -61441.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1444.4445
$ws.Range("I22").Value = 1444.4445
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1444.4445
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1149.4445
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1444.4445
$ws.Range("I27").Value = 1444.4445
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1444.4445
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1337.4445
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 3199.4
$ws.Range("I40").Value = 3199.4
$ws.Range("K40").Value = 3199.4
$ws.Range("M40").Value = -3063.4
$ws.Range("H46").Value = 463.63635
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 500
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -876
$ws.Range("H68").Value = 2278.9473
$ws.Range("I68").Value = 2154.5454
$ws.Range("J68").Value = 2450
$ws.Range("K68").Value = 2154.5454
$ws.Range("L68").Value = 2450
$ws.Range("M68").Value = -1405.5454
$ws.Range("N68").Value = -3948
$ws.Range("H71").Value = 2278.9473
$ws.Range("I71").Value = 2154.5454
$ws.Range("J71").Value = 2450
$ws.Range("K71").Value = 10772.727
$ws.Range("L71").Value = 12250
$ws.Range("M71").Value = -7028.726999999999
$ws.Range("N71").Value = -19738
$ws.Range("H135").Value = 48923.453
$ws.Range("J135").Value = 48923.453
$ws.Range("L135").Value = 48923.453
$ws.Range("N135").Value = -59063.453
$ws.Range("H136").Value = 5337.2
$ws.Range("I136").Value = 1697.4286
$ws.Range("J136").Value = 13830
$ws.Range("K136").Value = 5092.2858
$ws.Range("L136").Value = 41490
$ws.Range("M136").Value = -2542.2858
$ws.Range("N136").Value = -46590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 46512
$ws.Range("J131").Value = 46512
$ws.Range("L131").Value = 46512
$ws.Range("N131").Value = -56592
